# Daily attendance processing - 2026-01-19 12:01:19
# Normalize "Recorded By" (column G) entries so that "System" is listed
# first when a session shows it was recorded by both the instructor and
# the System (auto-recording), e.g. "dnasr281@gmail.com, System" ->
# "System, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}

Write-Host "Attendance recorder fields normalized."
